$p = $ppt.ActivePresentation

# --- Slide 2: title "TODOs von letztem Mal" -> "Grober Ablauf" ---
$s2 = $p.Slides.Item(2)
$titleShape = $null
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $shp = $s2.Shapes.Item($i)
    if ($shp.Name -eq "Titel 5") {
        $titleShape = $shp
    }
}
$titleShape.TextFrame.TextRange.Text = "Grober Ablauf"

# --- Slide 3: "Installation " + "von " runs merge into "Installation von " ---
$s3 = $p.Slides.Item(3)
$bodyShape = $null
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $shp = $s3.Shapes.Item($i)
    if ($shp.Name -eq "Textplatzhalter 1") {
        $bodyShape = $shp
    }
}
$bodyTextRange = $bodyShape.TextFrame.TextRange
for ($i = 1; $i -le $bodyTextRange.Paragraphs().Count; $i++) {
    $para = $bodyTextRange.Paragraphs($i)
    if ($para.Text.Trim() -eq "Installation von DevStack") {
        $para.Characters(1, 17).Text = "Installation von "
    }
}
